# Apply the "total for marking adjusted to reflect actual marks" edit:
#  1. "Code (50 marks)" -> "Code (25 marks)", where the text is split into
#     two runs with the _GoBack bookmark re-inserted between them (this is
#     what Word does when the last edit position is recorded mid-paragraph).
#  2. Remove the old _GoBack bookmark that used to sit after the "10" in the
#     "Data dictionary" marks cell (it moved to the edit above).

$d = $word.ActiveDocument

# --- 1. "Code (50 marks)" -> "Code (25" + bookmark + " marks)" ---------------
$r1 = $d.Content
$found1 = $r1.Find.Execute("Code (50 marks)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find 'Code (50 marks)' text"
}

$pPr1 = '<w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="2E75B5"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>'
$runProps1 = '<w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="2E75B5"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>'
$body1 = $pPr1 + '<w:r>' + $runProps1 + '<w:t>Code (25</w:t></w:r>' + '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + '<w:r>' + $runProps1 + '<w:t xml:space="preserve"> marks)</w:t></w:r>'

$xmlFrag1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="24B2BD9C" w14:textId="77777777" w:rsidR="00C908A3" w:rsidRDefault="00AE7983" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + $body1 + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$r1.InsertXML($xmlFrag1)

# --- 2. Remove the stale _GoBack bookmark after "10" -------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute("0 - 10", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find '0 - 10' text"
}

$pPr2 = '<w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>'
$runProps2 = '<w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>'
$body2 = $pPr2 + '<w:r>' + $runProps2 + '<w:t xml:space="preserve">0 - </w:t></w:r>' + '<w:r w:rsidR="006466B3">' + $runProps2 + '<w:t>10</w:t></w:r>'

$xmlFrag2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="090174DA" w14:textId="6941E19D" w:rsidR="00C908A3" w:rsidRDefault="00AE7983" w:rsidP="006466B3" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + $body2 + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$r2.InsertXML($xmlFrag2)

Write-Host "Edit complete"
